# Rename the DrawingML "name" metadata on the Pearson (footer) and BTec
# (header) logo pictures, swapping image1.png/image2.png and
# image1.jpg/image2.jpg the way the canonical edit did:
#
#   headers (BTec_Logo-Orange, wp:docPr/@name): image1.jpg -> image2.jpg
#   footers (PearsonLogo,      wp:docPr/@name): image2.png -> image1.png
#
# wdHeaderFooterPrimary / wdHeaderFooterFirstPage Word constants.
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Headers: BTec logo, image1.jpg -> image2.jpg ---------------------
    foreach ($kind in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $hdr = $sec.Headers.Item($kind)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    # InlineShape exposes no .Name setter of its own in the
                    # Word object model, but assigning it routes straight
                    # through to the underlying drawing's docPr/@name, so a
                    # direct assignment is enough here.
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # --- Footers: Pearson logo, image2.png -> image1.png ------------------
    foreach ($kind in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $ftr = $sec.Footers.Item($kind)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    # A plain InlineShape.Name assignment is silently
                    # ignored for footer pictures, so round-trip through a
                    # floating Shape (which does support .Name) and convert
                    # back to an inline picture afterwards.
                    $floatShp = $shp.ConvertToShape()
                    $floatShp.Name = "image1.png"
                    [void]$floatShp.ConvertToInlineShape()
                }
            }
        }
    }
}
